# Updated symbol list with latest crypto price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.28%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'41.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.93%"
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'1.82%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07674"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.55%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.627"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.24%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.9146"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'1.12%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'2.444"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.33%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1215"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'9.39%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1820"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2.40%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09106"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-0.58%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04216"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'1.48%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.1052"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.20%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001257"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.33%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005773"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.64%"
$ws.Range("E15").Style = "Normal"
$ws.Range("E17").Value = "'-0.32%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'1.26%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D20").Value = "'7.343"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'12.20%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'1.27%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2713"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-1.75%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04016"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-1.07%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001263"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'2.79%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004264"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'3.95%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001301"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'0.04%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02491"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'2.97%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05318"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'2.41%"
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'0.77%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1313"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'0.66%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.006502"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-7.75%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.001882"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-3.55%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008245"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-6.19%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3347"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'0.46%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006733"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-3.24%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'0.07%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.3719"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'1,100.23%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.003104"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-26.11%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'0.07%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.07%"
$ws.Range("E51").Style = "Normal"
